$d = $word.ActiveDocument

$d.Content.Find.Execute("57×99=5643", $true, $false, $false, $false, $false, $true, 1, $false, "74×96=7104", 2)
$d.Content.Find.Execute("25×49=1225", $true, $false, $false, $false, $false, $true, 1, $false, "61×88=5368", 2)
$d.Content.Find.Execute("14×49=686", $true, $false, $false, $false, $false, $true, 1, $false, "26×47=1222", 2)
$d.Content.Find.Execute("43×89=3827", $true, $false, $false, $false, $false, $true, 1, $false, "92×96=8832", 2)
$d.Content.Find.Execute("43×61=2623", $true, $false, $false, $false, $false, $true, 1, $false, "43×77=3311", 2)
$d.Content.Find.Execute("37×24=888", $true, $false, $false, $false, $false, $true, 1, $false, "98×13=1274", 2)
$d.Content.Find.Execute("36×66=2376", $true, $false, $false, $false, $false, $true, 1, $false, "69×23=1587", 2)
$d.Content.Find.Execute("13×71=923", $true, $false, $false, $false, $false, $true, 1, $false, "44×80=3520", 2)
$d.Content.Find.Execute("84×29=2436", $true, $false, $false, $false, $false, $true, 1, $false, "30×80=2400", 2)
$d.Content.Find.Execute("75×12=900", $true, $false, $false, $false, $false, $true, 1, $false, "17×71=1207", 2)
$d.Content.Find.Execute("31×41=1271", $true, $false, $false, $false, $false, $true, 1, $false, "58×71=4118", 2)
$d.Content.Find.Execute("84×64=5376", $true, $false, $false, $false, $false, $true, 1, $false, "91×69=6279", 2)
$d.Content.Find.Execute("52×59=3068", $true, $false, $false, $false, $false, $true, 1, $false, "62×51=3162", 2)
$d.Content.Find.Execute("95×60=5700", $true, $false, $false, $false, $false, $true, 1, $false, "15×90=1350", 2)
$d.Content.Find.Execute("48×62=2976", $true, $false, $false, $false, $false, $true, 1, $false, "42×54=2268", 2)
$d.Content.Find.Execute("36×46=1656", $true, $false, $false, $false, $false, $true, 1, $false, "71×84=5964", 2)
$d.Content.Find.Execute("71×77=5467", $true, $false, $false, $false, $false, $true, 1, $false, "86×84=7224", 2)
$d.Content.Find.Execute("12×21=252", $true, $false, $false, $false, $false, $true, 1, $false, "59×80=4720", 2)
$d.Content.Find.Execute("47×42=1974", $true, $false, $false, $false, $false, $true, 1, $false, "87×28=2436", 2)
$d.Content.Find.Execute("58×85=4930", $true, $false, $false, $false, $false, $true, 1, $false, "47×45=2115", 2)
$d.Content.Find.Execute("87×14=1218", $true, $false, $false, $false, $false, $true, 1, $false, "21×30=630", 2)
$d.Content.Find.Execute("44×84=3696", $true, $false, $false, $false, $false, $true, 1, $false, "84×95=7980", 2)
$d.Content.Find.Execute("96×40=3840", $true, $false, $false, $false, $false, $true, 1, $false, "14×66=924", 2)
$d.Content.Find.Execute("30×14=420", $true, $false, $false, $false, $false, $true, 1, $false, "55×21=1155", 2)
$d.Content.Find.Execute("79×18=1422", $true, $false, $false, $false, $false, $true, 1, $false, "80×17=1360", 2)
